$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix NACA 0012 Re (1) entries (rows 18-19) ---
$ws.Range("D18").Value = "2e6 - 4e6"
$ws.Range("D19").Value = "1e6 - 2e6"

# --- Row 20 (CAST 7): remove bold, add a note ---
$ws.Range("A20:N20").Font.Bold = $false
$ws.Range("N20").Value = "1. Fixed transition"

# --- Row 33: repurpose old SSC-B08 placeholder row into SSC-A08 ---
$ws.Range("A33").Value = "SSC-A08"
$ws.Range("B33").Value = "Need to do"

# --- Row 34: add new SKF 1.1 entry ---
$ws.Range("A34").Value = "SKF 1.1"
$ws.Range("B34").Value = "2.5, 5.0"
$ws.Range("C34").Value = 0.76
$ws.Range("D34").Value = "2.2e6, 7e6"
$ws.Range("E34").Value = 0.01
$ws.Range("G34").Value = 2
$ws.Range("H34").Value = "AGARD AR 138"
$ws.Range("K34").Value = "Wing"
$ws.Range("L34").Value = "Supercritical"
$ws.Range("M34").Value = "Cambered"
$ws.Range("N34").Value = "1. Experiment is w/ Maneuver flap, only baseline is digitized"

# --- Resize the table to include the two new rows ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:N34"))
